# Auto-generated edit script: update cached market-data columns (H-N)
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 5: refreshed market-price snapshot
$ws.Range("H5").Value = 96
$ws.Range("I5").Value = 96
$ws.Range("K5").Value = 96
$ws.Range("M5").Value = 19
# ALC row 39: refreshed market-price snapshot
$ws.Range("H39").Value = 674.1818
$ws.Range("I39").Value = 739.6
$ws.Range("J39").Value = 20
$ws.Range("K39").Value = 2218.8
$ws.Range("L39").Value = 60
$ws.Range("M39").Value = -1922.8
$ws.Range("N39").Value = -652
# ALC row 113: refreshed market-price snapshot
$ws.Range("H113").Value = 4239.75
$ws.Range("I113").Value = 3979.5
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 3979.5
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -725.5
$ws.Range("N113").Value = -11008
# ALC row 132: refreshed market-price snapshot
$ws.Range("H132").Value = 1378.3334
$ws.Range("I132").Value = 1378.3334
$ws.Range("K132").Value = 4135.0002
$ws.Range("M132").Value = -1605.0002

$ws = $wb.Worksheets.Item("ARM")
# ARM row 19: refreshed market-price snapshot
$ws.Range("H19").Value = 1090
$ws.Range("I19").Value = 1090
$ws.Range("K19").Value = 1090
$ws.Range("M19").Value = -861
# ARM row 32: refreshed market-price snapshot
$ws.Range("H32").Value = 7043.8
$ws.Range("I32").Value = 7442.846
$ws.Range("K32").Value = 7442.846
$ws.Range("M32").Value = -7155.846
# ARM row 61: refreshed market-price snapshot
$ws.Range("H61").Value = 4545.0557
$ws.Range("I61").Value = 1585
$ws.Range("J61").Value = 5390.7856
$ws.Range("K61").Value = 1585
$ws.Range("L61").Value = 5390.7856
$ws.Range("M61").Value = -1373
$ws.Range("N61").Value = -5814.7856
# ARM row 136: refreshed market-price snapshot
$ws.Range("H136").Value = 4545.0557
$ws.Range("I136").Value = 1585
$ws.Range("J136").Value = 5390.7856
$ws.Range("K136").Value = 4755
$ws.Range("L136").Value = 16172.3568
$ws.Range("M136").Value = -2205
$ws.Range("N136").Value = -21272.3568

$ws = $wb.Worksheets.Item("BSM")
# BSM row 19: refreshed market-price snapshot
$ws.Range("H19").Value = 3000
$ws.Range("J19").Value = 3000
$ws.Range("L19").Value = 3000
$ws.Range("N19").Value = -3346
# BSM row 81: refreshed market-price snapshot
$ws.Range("H81").Value = 38962
$ws.Range("J81").Value = 38962
$ws.Range("L81").Value = 38962
$ws.Range("N81").Value = -41084
# BSM row 84: refreshed market-price snapshot
$ws.Range("H84").Value = 38962
$ws.Range("J84").Value = 38962
$ws.Range("L84").Value = 116886
$ws.Range("N84").Value = -127494
# BSM row 86: refreshed market-price snapshot
$ws.Range("H86").Value = 450
$ws.Range("I86").Value = 450
$ws.Range("K86").Value = 450
$ws.Range("M86").Value = 673
# BSM row 89: refreshed market-price snapshot
$ws.Range("H89").Value = 450
$ws.Range("I89").Value = 450
$ws.Range("K89").Value = 2250
$ws.Range("M89").Value = 3366

$ws = $wb.Worksheets.Item("CRP")
# CRP row 56: refreshed market-price snapshot
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
# CRP row 58: refreshed market-price snapshot
$ws.Range("H58").Value = 1688.55
$ws.Range("J58").Value = 2198.5
$ws.Range("L58").Value = 2198.5
$ws.Range("N58").Value = -2604.5
# CRP row 99: refreshed market-price snapshot
$ws.Range("H99").Value = 2886.2
$ws.Range("I99").Value = 2649
$ws.Range("J99").Value = 3242
$ws.Range("K99").Value = 2649
$ws.Range("L99").Value = 3242
$ws.Range("M99").Value = -1151
$ws.Range("N99").Value = -6238
# CRP row 126: refreshed market-price snapshot
$ws.Range("H126").Value = 2886.2
$ws.Range("I126").Value = 2649
$ws.Range("J126").Value = 3242
$ws.Range("K126").Value = 7947
$ws.Range("L126").Value = 9726
$ws.Range("M126").Value = -5477
$ws.Range("N126").Value = -14666
# CRP row 132: refreshed market-price snapshot
$ws.Range("H132").Value = 3089.5
$ws.Range("I132").Value = 3981.3333
$ws.Range("J132").Value = 2197.6667
$ws.Range("K132").Value = 11943.9999
$ws.Range("L132").Value = 6593.000100000001
$ws.Range("M132").Value = -9413.999899999999
$ws.Range("N132").Value = -11653.0001
# CRP row 136: refreshed market-price snapshot
$ws.Range("H136").Value = 1688.55
$ws.Range("J136").Value = 2198.5
$ws.Range("L136").Value = 6595.5
$ws.Range("N136").Value = -11695.5
# CRP row 141: refreshed market-price snapshot
$ws.Range("H141").Value = 1000000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 1000000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 1000000
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -1010360

$ws = $wb.Worksheets.Item("CUL")
# CUL row 4: refreshed market-price snapshot
$ws.Range("H4").Value = 12222913
$ws.Range("I4").Value = 12222913
$ws.Range("K4").Value = 36668739
$ws.Range("M4").Value = -36668627
# CUL row 31: refreshed market-price snapshot
$ws.Range("H31").Value = 250
$ws.Range("I31").Value = 250
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 750
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -462
$ws.Range("N31").ClearContents()
# CUL row 123: refreshed market-price snapshot
$ws.Range("H123").Value = 8996
$ws.Range("I123").Value = 7326
$ws.Range("K123").Value = 21978
$ws.Range("M123").Value = -19528

$ws = $wb.Worksheets.Item("GSM")
# GSM row 70: refreshed market-price snapshot
$ws.Range("H70").Value = 8647.049999999999
$ws.Range("I70").Value = 7006.3335
$ws.Range("J70").Value = 9989.454
$ws.Range("K70").Value = 7006.3335
$ws.Range("L70").Value = 9989.454
$ws.Range("M70").Value = -6736.3335
$ws.Range("N70").Value = -10529.454
# GSM row 73: refreshed market-price snapshot
$ws.Range("H73").Value = 8647.049999999999
$ws.Range("I73").Value = 7006.3335
$ws.Range("J73").Value = 9989.454
$ws.Range("K73").Value = 7006.3335
$ws.Range("L73").Value = 9989.454
$ws.Range("M73").Value = -6070.3335
$ws.Range("N73").Value = -11861.454
# GSM row 132: refreshed market-price snapshot
$ws.Range("H132").Value = 1648.6666
$ws.Range("I132").Value = 1474
$ws.Range("J132").Value = 1998
$ws.Range("K132").Value = 4422
$ws.Range("L132").Value = 5994
$ws.Range("M132").Value = -1892
$ws.Range("N132").Value = -11054

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22: refreshed market-price snapshot
$ws.Range("H22").Value = 4201.2856
$ws.Range("I22").Value = 4165.2
$ws.Range("J22").Value = 4291.5
$ws.Range("K22").Value = 4165.2
$ws.Range("L22").Value = 4291.5
$ws.Range("M22").Value = -3870.2
$ws.Range("N22").Value = -4881.5
# LTW row 27: refreshed market-price snapshot
$ws.Range("H27").Value = 4201.2856
$ws.Range("I27").Value = 4165.2
$ws.Range("J27").Value = 4291.5
$ws.Range("K27").Value = 4165.2
$ws.Range("L27").Value = 4291.5
$ws.Range("M27").Value = -4058.2
$ws.Range("N27").Value = -4505.5
# LTW row 45: refreshed market-price snapshot
$ws.Range("H45").Value = 29900
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 29900
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 29900
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -30714
# LTW row 108: refreshed market-price snapshot
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680
# LTW row 121: refreshed market-price snapshot
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
# LTW row 122: refreshed market-price snapshot
$ws.Range("H122").Value = 6893.7646
$ws.Range("I122").Value = 6540.727
$ws.Range("K122").Value = 19622.181
$ws.Range("M122").Value = -17172.181
# LTW row 133: refreshed market-price snapshot
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
# LTW row 134: refreshed market-price snapshot
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
# LTW row 135: refreshed market-price snapshot
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# WVR row 100: refreshed market-price snapshot
$ws.Range("H100").Value = 9092245
$ws.Range("I100").Value = 10001401
$ws.Range("K100").Value = 20002802
$ws.Range("M100").Value = -20002261
# WVR row 132: refreshed market-price snapshot
$ws.Range("H132").Value = 8488.799999999999
$ws.Range("I132").Value = 7765.4443
$ws.Range("J132").Value = 14999
$ws.Range("K132").Value = 23296.3329
$ws.Range("L132").Value = 44997
$ws.Range("M132").Value = -20766.3329
$ws.Range("N132").Value = -50057
# WVR row 136: refreshed market-price snapshot
$ws.Range("H136").Value = 4090.375
$ws.Range("I136").Value = 4128.722
$ws.Range("J136").Value = 3975.3333
$ws.Range("K136").Value = 12386.166
$ws.Range("L136").Value = 11925.9999
$ws.Range("M136").Value = -9836.165999999999
$ws.Range("N136").Value = -17025.9999
